$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Signatory change: the NKO representative moves from the Chairman of the
# Board (acting on the charter) to a Deputy Chairman (acting by power of
# attorney).
$ws.Range("A41").Value = "Заместитель Председателя Правления"
$ws.Range("A42").Value = "Лебедева Л.В."
$ws.Range("A44").Value = "Доверенность № 13 от 21.02.2020"
